$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value2 = "2026-02-16 23:48:28"
$ws.Cells.Item(3, 5).Value2 = "2026-02-16 23:48:30"
$ws.Cells.Item(3, 14).Value2 = "-5.8 °C 23:29 TU"
$ws.Cells.Item(3, 15).Value2 = "-1.4 °C"
$ws.Cells.Item(4, 5).Value2 = "2026-02-16 23:48:33"
$ws.Cells.Item(4, 14).Value2 = "6.6 °C 23:26 TU"
$ws.Cells.Item(4, 15).Value2 = "12.9 °C"
$ws.Cells.Item(5, 5).Value2 = "2026-02-16 23:48:35"
$ws.Cells.Item(5, 9).Value2 = "25.2 mm"
$ws.Cells.Item(5, 14).Value2 = "-6.0 °C 23:29 TU"
$ws.Cells.Item(5, 15).Value2 = "-1.2 °C"
$ws.Cells.Item(6, 5).Value2 = "2026-02-16 23:48:38"
$ws.Cells.Item(7, 5).Value2 = "2026-02-16 23:48:40"
$ws.Cells.Item(7, 10).Value2 = "1013.3 hPa"
$ws.Cells.Item(8, 5).Value2 = "2026-02-16 23:48:43"
$ws.Cells.Item(8, 8).NumberFormat = "@"
$ws.Cells.Item(8, 8).Value2 = "67%"
$ws.Cells.Item(9, 5).Value2 = "2026-02-16 23:48:45"
$ws.Cells.Item(9, 8).NumberFormat = "@"
$ws.Cells.Item(9, 8).Value2 = "67%"
$ws.Cells.Item(9, 12).Value2 = "52.9 km/h - 347º 23:24 TU"
$ws.Cells.Item(9, 15).Value2 = "11.7 °C"
$ws.Cells.Item(10, 5).Value2 = "2026-02-16 23:48:48"
$ws.Cells.Item(11, 5).Value2 = "2026-02-16 23:48:49"
$ws.Cells.Item(11, 8).NumberFormat = "@"
$ws.Cells.Item(11, 8).Value2 = "75%"
$ws.Cells.Item(11, 15).Value2 = "7.2 °C"
$ws.Cells.Item(12, 5).Value2 = "2026-02-16 23:48:50"
$ws.Cells.Item(12, 8).NumberFormat = "@"
$ws.Cells.Item(12, 8).Value2 = "75%"
$ws.Cells.Item(13, 5).Value2 = "2026-02-16 23:48:51"
$ws.Cells.Item(14, 5).Value2 = "2026-02-16 23:48:52"
$ws.Cells.Item(14, 15).Value2 = "16.1 °C"
$ws.Cells.Item(15, 5).Value2 = "2026-02-16 23:48:53"
$ws.Cells.Item(15, 8).NumberFormat = "@"
$ws.Cells.Item(15, 8).Value2 = "63%"
$ws.Cells.Item(15, 15).Value2 = "12.0 °C"
$ws.Cells.Item(16, 5).Value2 = "2026-02-16 23:48:54"
$ws.Cells.Item(16, 8).NumberFormat = "@"
$ws.Cells.Item(16, 8).Value2 = "78%"
$ws.Cells.Item(16, 14).Value2 = "-5.5 °C 23:11 TU"
$ws.Cells.Item(16, 15).Value2 = "-0.6 °C"
$ws.Cells.Item(17, 5).Value2 = "2026-02-16 23:48:55"
$ws.Cells.Item(17, 14).Value2 = "2.8 °C 23:21 TU"
$ws.Cells.Item(18, 5).Value2 = "2026-02-16 23:48:57"
$ws.Cells.Item(18, 8).NumberFormat = "@"
$ws.Cells.Item(18, 8).Value2 = "77%"
$ws.Cells.Item(18, 15).Value2 = "10.4 °C"
$ws.Cells.Item(19, 5).Value2 = "2026-02-16 23:48:58"
$ws.Cells.Item(20, 5).Value2 = "2026-02-16 23:48:59"
$ws.Cells.Item(20, 8).NumberFormat = "@"
$ws.Cells.Item(20, 8).Value2 = "97%"
$ws.Cells.Item(20, 12).Value2 = "74.5 km/h - 339º 23:24 TU"
$ws.Cells.Item(21, 5).Value2 = "2026-02-16 23:49:00"
$ws.Cells.Item(22, 5).Value2 = "2026-02-16 23:49:02"
$ws.Cells.Item(23, 5).Value2 = "2026-02-16 23:49:05"
$ws.Cells.Item(23, 9).Value2 = "16.3 mm"
$ws.Cells.Item(23, 14).Value2 = "-5.2 °C 23:19 TU"
$ws.Cells.Item(23, 15).Value2 = "-1.2 °C"
$ws.Cells.Item(24, 5).Value2 = "2026-02-16 23:49:07"
$ws.Cells.Item(25, 5).Value2 = "2026-02-16 23:49:10"
$ws.Cells.Item(25, 8).NumberFormat = "@"
$ws.Cells.Item(25, 8).Value2 = "84%"
$ws.Cells.Item(26, 5).Value2 = "2026-02-16 23:49:12"
$ws.Cells.Item(27, 5).Value2 = "2026-02-16 23:49:14"
$ws.Cells.Item(27, 8).NumberFormat = "@"
$ws.Cells.Item(27, 8).Value2 = "80%"
$ws.Cells.Item(27, 14).Value2 = "-1.4 °C 23:23 TU"
$ws.Cells.Item(28, 5).Value2 = "2026-02-16 23:49:17"
$ws.Cells.Item(28, 8).NumberFormat = "@"
$ws.Cells.Item(28, 8).Value2 = "74%"
$ws.Cells.Item(28, 15).Value2 = "9.3 °C"
$ws.Cells.Item(29, 5).Value2 = "2026-02-16 23:49:19"
$ws.Cells.Item(29, 8).NumberFormat = "@"
$ws.Cells.Item(29, 8).Value2 = "79%"
$ws.Cells.Item(29, 15).Value2 = "11.0 °C"
$ws.Cells.Item(30, 5).Value2 = "2026-02-16 23:49:21"
$ws.Cells.Item(30, 8).NumberFormat = "@"
$ws.Cells.Item(30, 8).Value2 = "67%"
$ws.Cells.Item(31, 5).Value2 = "2026-02-16 23:49:24"
$ws.Cells.Item(31, 10).Value2 = "1011.8 hPa"
$ws.Cells.Item(31, 15).Value2 = "14.1 °C"
$ws.Cells.Item(32, 5).Value2 = "2026-02-16 23:49:26"
$ws.Cells.Item(32, 15).Value2 = "8.6 °C"
$ws.Cells.Item(33, 5).Value2 = "2026-02-16 23:49:29"
$ws.Cells.Item(33, 8).NumberFormat = "@"
$ws.Cells.Item(33, 8).Value2 = "71%"
$ws.Cells.Item(34, 5).Value2 = "2026-02-16 23:49:31"
$ws.Cells.Item(34, 15).Value2 = "3.1 °C"
$ws.Cells.Item(35, 5).Value2 = "2026-02-16 23:49:34"
$ws.Cells.Item(35, 8).NumberFormat = "@"
$ws.Cells.Item(35, 8).Value2 = "77%"
$ws.Cells.Item(35, 9).Value2 = "3.9 mm"
$ws.Cells.Item(35, 10).Value2 = "1016.6 hPa"
$ws.Cells.Item(35, 12).Value2 = "82.1 km/h - 247º 23:22 TU"
$ws.Cells.Item(35, 15).Value2 = "9.3 °C"
$ws.Cells.Item(36, 5).Value2 = "2026-02-16 23:49:36"
$ws.Cells.Item(36, 8).NumberFormat = "@"
$ws.Cells.Item(36, 8).Value2 = "68%"
$ws.Cells.Item(37, 5).Value2 = "2026-02-16 23:49:39"
$ws.Cells.Item(37, 8).NumberFormat = "@"
$ws.Cells.Item(37, 8).Value2 = "82%"
$ws.Cells.Item(37, 12).Value2 = "42.1 km/h - 252º 23:25 TU"
$ws.Cells.Item(38, 5).Value2 = "2026-02-16 23:49:41"
$ws.Cells.Item(38, 8).NumberFormat = "@"
$ws.Cells.Item(38, 8).Value2 = "72%"
$ws.Cells.Item(39, 5).Value2 = "2026-02-16 23:49:44"
$ws.Cells.Item(39, 12).Value2 = "74.2 km/h - 225º 23:26 TU"
$ws.Cells.Item(39, 15).Value2 = "-0.3 °C"
$ws.Cells.Item(40, 5).Value2 = "2026-02-16 23:49:46"
$ws.Cells.Item(40, 15).Value2 = "6.9 °C"
$ws.Cells.Item(41, 5).Value2 = "2026-02-16 23:49:48"
$ws.Cells.Item(42, 5).Value2 = "2026-02-16 23:49:51"
$ws.Cells.Item(43, 5).Value2 = "2026-02-16 23:49:53"
$ws.Cells.Item(44, 5).Value2 = "2026-02-16 23:49:56"
$ws.Cells.Item(44, 14).Value2 = "-4.9 °C 23:19 TU"
$ws.Cells.Item(45, 5).Value2 = "2026-02-16 23:49:58"
$ws.Cells.Item(45, 10).Value2 = "1018.1 hPa"
$ws.Cells.Item(45, 14).Value2 = "2.1 °C 23:28 TU"
$ws.Cells.Item(45, 15).Value2 = "4.7 °C"
$ws.Cells.Item(46, 5).Value2 = "2026-02-16 23:50:01"
$ws.Cells.Item(46, 8).NumberFormat = "@"
$ws.Cells.Item(46, 8).Value2 = "57%"
